{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (plus the blank paragraph that trails it), left over from the scraped\n// site chrome, right after the bibliography entry for FLEMMING/GON\u00c7ALVES.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst COPYRIGHT_MARK = \"Contact: luizeleno@usp.br\";\nconst JUPITER_MARK = \"Ver no Jupiter\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(COPYRIGHT_MARK) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const prev = target.getPrevious(); // \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n  const next = target.getNext(); // blank paragraph right before the page break\n\n  prev.load(\"text\");\n  next.load(\"text\");\n  await context.sync();\n\n  // Sanity-check before deleting, so we never nuke unrelated content if the\n  // document shape differs from what we expect.\n  if (prev.text && prev.text.indexOf(JUPITER_MARK) !== -1) {\n    prev.delete();\n  }\n  target.delete();\n  if (next.text === \"\") {\n    next.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (plus the blank paragraph that trails it), left over from the scraped\n# site chrome, right after the bibliography entry for FLEMMING/GON\u00c7ALVES.\n\n$d = $word.ActiveDocument\n\n$copyrightMark = \"Contact: luizeleno@usp.br\"\n$jupiterMark = \"Ver no Jupiter\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$copyrightMark*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $prevP = $target.Previous()\n    $nextP = $target.Next()\n\n    $startPos = $target.Range.Start\n    $endPos = $target.Range.End\n\n    if (($prevP -ne $null) -and ($prevP.Range.Text -like \"*$jupiterMark*\")) {\n        $startPos = $prevP.Range.Start\n    }\n    if (($nextP -ne $null) -and ($nextP.Range.Text.Trim() -eq \"\")) {\n        $endPos = $nextP.Range.End\n    }\n\n    $delRange = $d.Range($startPos, $endPos)\n    $delRange.Delete()\n}\n"}
